$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing "2026/12/29" block (currently rows 776-777),
# shifting all subsequent rows down by two. This grows the data from A1:D817 to A1:D819.
$ws.Range("A776:D777").Insert()

# Force column A on the new rows to be stored as plain text (matching the rest of the
# date column) instead of being auto-converted to a date serial number.
$ws.Range("A776:A777").NumberFormat = "@"

# Row 776: 2026/02/05, 木, 18, 201
$ws.Range("A776").Value = "2026/02/05"
$ws.Range("B776").Value = "木"
$ws.Range("C776").Value = 18
$ws.Range("D776").Value = 201

# Row 777: 2026/02/05, 木, 22, 201
$ws.Range("A777").Value = "2026/02/05"
$ws.Range("B777").Value = "木"
$ws.Range("C777").Value = 22
$ws.Range("D777").Value = 201

# Drop the temporary text formatting so the new cells carry no explicit style,
# matching the unstyled data cells elsewhere in the sheet.
$ws.Range("A776:D777").ClearFormats()
